# Update Croatia configuration calibrated input variables
# Sheet "strategy_id-6008" row 2 contains a decay curve of values; update
# the segment from V2 through AR2 with the recalculated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-6008")

$values = @{
    "V2"  = 0.88
    "W2"  = 0.8539130434782609
    "X2"  = 0.8286956521739131
    "Y2"  = 0.8043478260869564
    "Z2"  = 0.7808695652173914
    "AA2" = 0.7582608695652174
    "AB2" = 0.7365217391304347
    "AC2" = 0.7156521739130435
    "AD2" = 0.6956521739130435
    "AE2" = 0.6765217391304349
    "AF2" = 0.6582608695652175
    "AG2" = 0.6408695652173912
    "AH2" = 0.6243478260869566
    "AI2" = 0.6086956521739131
    "AJ2" = 0.5939130434782609
    "AK2" = 0.5800000000000001
    "AL2" = 0.5669565217391305
    "AM2" = 0.5547826086956522
    "AN2" = 0.5434782608695652
    "AO2" = 0.5347826086956522
    "AP2" = 0.5260869565217391
    "AQ2" = 0.5173913043478261
    "AR2" = 0.508695652173913
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
